# "actu et nommay deplace"
#
# The Nommay event (row 49, "Gentlemen de Nommay") is rescheduled from
# "Dim 20 Septembre" to "Dim 18 Octobre" and moves down the calendar to
# sit in date order between the Giromagny cyclo-cross (Dim 11 Octobre,
# row 55) and the Heimsbrunn cyclo-cross (Sam 24 Octobre, row 56). Every
# row that used to sit between the old and new Nommay slot (rows 50-55)
# shifts up by one row to fill the gap left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param([int]$Row, [string]$A, [string]$B, [string]$C, [string]$D, [string]$E, [string]$G)
    $ws.Range("A$Row").Value = $A
    $ws.Range("B$Row").Value = $B
    $ws.Range("C$Row").Value = $C
    $ws.Range("D$Row").Value = $D
    $ws.Range("E$Row").Value = $E
    $ws.Range("G$Row").Value = $G
}

# old row 50 -> row 49
Set-Row 49 "Sam 26 Septembre" "5e VTT MS Automobile Rixheim  " "SSOL Habsheim" "VTT" "rixheim" "Annulé (refus de l'ONF)"

# old row 51 -> row 50
Set-Row 50 "Dim 27 Septembre" "La Geko-Bikes à Didenheim" "UC Lutterbach VTT" "VTT" "didenheim" ""

# old row 52 -> row 51
Set-Row 51 "Sam 3 Octobre" "La Grimpée du Salbert Trophée Maurice Voirol à Cravanche  " "ACT Belfort" "Grimpée" "salbert" ""

# old row 53 -> row 52
Set-Row 52 "Dim 4 Octobre" "Prix Technochape à Retzwiller" "ASCL Montreux-Vieux" "VTT" "retzwillerbis" "Championnat d'Alsace"

# old row 54 -> row 53
Set-Row 53 "Sam 10 Octobre" "Cyclo-cross de Morschwiller le bas" "SOS Lutterbach" "Cyclo-cross" "morschwiller" ""

# old row 55 -> row 54
Set-Row 54 "Dim 11 Octobre" "3e Cyclo-cross de Giromagny. Epreuve FFC ouverte aux FSGT" "US Giromagny VTT" "Cyclo-cross" "giromagny" ""

# old row 49 (Nommay) -> row 55, date updated to "Dim 18 Octobre"
Set-Row 55 "Dim 18 Octobre" "Gentlemen de Nommay (épreuve FFC ouverte aux FSGT)*" "CCI Nommay" "Route" "nommay" ""

# Update the view state: scroll so row 24 is near the top, and select
# B55 - the cell where the relocated Nommay row now lives.
$ws.Range("A24").Select() | Out-Null
$ws.Range("B55").Select() | Out-Null
